{"js": "// Apply the \"Knight's Keep\" copy refresh described in the commit\n// (\"Added many more features\"): retitle the headline/meta-title,\n// rewrite the pros/cons bullet list, and refresh the closing\n// bold/italic summary paragraphs.\n//\n// Each entry is searched for literally (exact, case-sensitive text)\n// and replaced in place via Range.insertText(..., \"Replace\"), which\n// keeps the owning run's formatting (bold/italic/heading style) intact.\nconst replacements = [\n  [\n    \"Play Knight's Keep for Free - Retro-Themed 3D Slot Game\",\n    \"Play Knight\\u2019s Keep Free: Retro-Themed Slot Game with Exciting Features\",\n  ],\n  [\n    \"192 paylines provide a higher chance of winning\",\n    \"Combines classic gameplay features with modern 3D graphics\",\n  ],\n  [\n    \"The black knight symbol acts as a Wild and a 10x multiplier\",\n    \"192 paylines for a higher chance of winning\",\n  ],\n  [\n    \"Impressive 3D graphics create a medieval world\",\n    \"Black knight symbol acts as a Wild and a 10x multiplier\",\n  ],\n  [\n    \"Free spins bonus feature can award up to 25 free spins\",\n    \"Free spins bonus feature that can award up to 25 free spins\",\n  ],\n  [\n    \"May not appeal to those who don't enjoy retro or medieval themes\",\n    \"Limited number of bonus features\",\n  ],\n  [\n    \"Not a unique game in terms of features\",\n    \"May not appeal to players who prefer non-medieval themes\",\n  ],\n  [\n    \"Experience Knight's Keep, a retro-themed slot game with 192 paylines, free spins, and exciting multiplier symbols. Play for free and discover a medieval world.\",\n    \"Experience the thrill of Knight\\u2019s Keep, a retro-themed slot game with modern 3D graphics. Play now for free!\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Knight's Keep\" copy refresh described in the commit\n# (\"Added many more features\"): retitle the headline/meta-title,\n# rewrite the pros/cons bullet list, and refresh the closing\n# bold/italic summary paragraphs.\n#\n# Each old -> new pair is run through Find/Replace (wdReplaceAll) on\n# the whole document Range, which finds every literal occurrence\n# (the title text appears twice: Heading1 + the bold closer) and\n# replaces it in place, preserving each run's formatting.\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$replacements = @(\n    @{\n        Find    = \"Play Knight's Keep for Free - Retro-Themed 3D Slot Game\"\n        Replace = \"Play Knight\" + [char]0x2019 + \"s Keep Free: Retro-Themed Slot Game with Exciting Features\"\n    },\n    @{\n        Find    = \"192 paylines provide a higher chance of winning\"\n        Replace = \"Combines classic gameplay features with modern 3D graphics\"\n    },\n    @{\n        Find    = \"The black knight symbol acts as a Wild and a 10x multiplier\"\n        Replace = \"192 paylines for a higher chance of winning\"\n    },\n    @{\n        Find    = \"Impressive 3D graphics create a medieval world\"\n        Replace = \"Black knight symbol acts as a Wild and a 10x multiplier\"\n    },\n    @{\n        Find    = \"Free spins bonus feature can award up to 25 free spins\"\n        Replace = \"Free spins bonus feature that can award up to 25 free spins\"\n    },\n    @{\n        Find    = \"May not appeal to those who don't enjoy retro or medieval themes\"\n        Replace = \"Limited number of bonus features\"\n    },\n    @{\n        Find    = \"Not a unique game in terms of features\"\n        Replace = \"May not appeal to players who prefer non-medieval themes\"\n    },\n    @{\n        Find    = \"Experience Knight's Keep, a retro-themed slot game with 192 paylines, free spins, and exciting multiplier symbols. Play for free and discover a medieval world.\"\n        Replace = \"Experience the thrill of Knight\" + [char]0x2019 + \"s Keep, a retro-themed slot game with modern 3D graphics. Play now for free!\"\n    }\n)\n\nforeach ($item in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($item.Find, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $item.Replace, $wdReplaceAll)\n}\n"}
